$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8; existing rows 8.. shift down by one
# (this mirrors Excel's Rows.Insert behaviour, including carrying down
# the formatting of the row above into the newly inserted row).
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new weekly record.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C8").Value = "Ñuble"
$ws.Range("D8").Value = 45063
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 100112043
$ws.Range("G8").Value = "Pepino dulce"
$ws.Range("H8").Value = "Cultivar IV Región"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 60
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 15000
$ws.Range("N8").Value = "$/bandeja 18 kilos"
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 833
$ws.Range("Q8").Value = 18
$ws.Range("R8").Value = "Hortaliza"
